$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 16,33
for ($r = 0; $r -lt 16; $r++) { for ($c = 0; $c -lt 33; $c++) { $data[$r,$c] = $null } }
$data[0,0] = "Name(ID)"
$data[0,1] = "sprite_name"
$data[0,2] = "chat_id"
$data[0,3] = "max_hp"
$data[0,4] = "max_shield"
$data[0,5] = "max_stagger"
$data[0,6] = "ATK"
$data[0,7] = "DEF"
$data[0,8] = "SPEED "
$data[0,9] = "ACC (%)"
$data[0,10] = "EVADE"
$data[0,11] = "AI type"
$data[0,12] = "AI parameters"
$data[0,13] = "null wk/res"
$data[0,14] = "agni wk/res"
$data[0,15] = "cryo wk/res"
$data[0,16] = "bolt wk/res"
$data[0,17] = "SpellGroup 1"
$data[0,18] = "Spell1 (root)"
$data[0,19] = "Spell1 (elem)"
$data[0,20] = "Spell1(style)"
$data[0,21] = "Spell x (root)"
$data[0,22] = "Spell x (elem)"
$data[0,23] = "Spell x (style)"
$data[0,24] = "Next Group or END"
$data[1,0] = "Slime"
$data[1,1] = "frog_mario"
$data[1,2] = "enemy_general_1"
$data[1,3] = 10
$data[1,4] = 0
$data[1,5] = 1
$data[1,6] = 1
$data[1,7] = 0
$data[1,8] = 0.5
$data[1,9] = 1.5
$data[1,10] = 0
$data[1,11] = "Attacker1"
$data[1,12] = "none"
$data[1,13] = 1
$data[1,14] = 0.5
$data[1,15] = 1.5
$data[1,16] = 1
$data[1,17] = "GROUP/DEFAULT"
$data[1,18] = "sword"
$data[1,19] = "null"
$data[1,20] = "null"
$data[1,21] = "hammer"
$data[1,22] = "null"
$data[1,23] = "null"
$data[1,24] = "GROUP/HEALTH_LOW"
$data[1,25] = "quake"
$data[1,26] = "null"
$data[1,27] = "null"
$data[1,28] = "END"
$data[2,0] = "The Evil Eye"
$data[2,1] = "frog_mario"
$data[2,2] = "enemy_general_1"
$data[2,3] = 50
$data[2,4] = 0
$data[2,5] = 2
$data[2,6] = 1.25
$data[2,7] = 0.2
$data[2,8] = 0.75
$data[2,9] = 1
$data[2,10] = 0
$data[2,11] = "HealthLow1"
$data[2,12] = "75,25"
$data[2,13] = -1
$data[2,14] = 2
$data[2,15] = 0.75
$data[2,16] = -2
$data[2,17] = "GROUP/DEFAULT"
$data[2,18] = "sword"
$data[2,19] = "null"
$data[2,20] = "null"
$data[2,21] = "lance"
$data[2,22] = "agni"
$data[2,23] = "null"
$data[2,24] = "GROUP/HEALTH_LOW"
$data[2,25] = "quake"
$data[2,26] = "null"
$data[2,27] = "null"
$data[2,28] = "END"
$data[3,0] = "Ladon"
$data[3,1] = "spr_bt_ladon"
$data[3,2] = "enemy_general_1"
$data[3,3] = 100
$data[3,4] = 0
$data[3,5] = 3
$data[3,6] = 0.95
$data[3,7] = -0.2
$data[3,8] = 0.3
$data[3,9] = 1
$data[3,10] = 10
$data[3,11] = "Attacker1"
$data[3,12] = "none"
$data[3,13] = 1.2
$data[3,14] = 1
$data[3,15] = 1
$data[3,16] = 1
$data[3,17] = "GROUP/DEFAULT"
$data[3,18] = "lance"
$data[3,19] = "null"
$data[3,20] = "null"
$data[3,21] = "lance"
$data[3,22] = "null"
$data[3,23] = "aimed"
$data[3,24] = "GROUP/HEALTH_LOW"
$data[3,25] = "quake"
$data[3,26] = "null"
$data[3,27] = "null"
$data[3,28] = "END"
$data[4,0] = "Lilim"
$data[4,1] = "spr_bt_lilim"
$data[4,2] = "enemy_general_1"
$data[4,3] = 50
$data[4,4] = 0
$data[4,5] = 2
$data[4,6] = 1.25
$data[4,7] = 0.2
$data[4,8] = 0.75
$data[4,9] = 1
$data[4,10] = 0
$data[4,11] = "HealthLow1"
$data[4,12] = "75,25"
$data[4,13] = -1
$data[4,14] = 2
$data[4,15] = 0.75
$data[4,16] = -2
$data[4,17] = "GROUP/DEFAULT"
$data[4,18] = "sword"
$data[4,19] = "null"
$data[4,20] = "null"
$data[4,21] = "lance"
$data[4,22] = "agni"
$data[4,23] = "null"
$data[4,24] = "GROUP/HEALTH_LOW"
$data[4,25] = "quake"
$data[4,26] = "null"
$data[4,27] = "null"
$data[4,28] = "END"
$data[5,0] = "Changeling"
$data[5,1] = "spr_bt_changeling_placeholder"
$data[5,2] = "changeling_1"
$data[5,3] = 10
$data[5,4] = 0
$data[5,5] = 1
$data[5,6] = 1
$data[5,7] = 0
$data[5,8] = 0.5
$data[5,9] = 1
$data[5,10] = 0
$data[5,11] = "Attacker1"
$data[5,12] = "none"
$data[5,13] = 1
$data[5,14] = 1
$data[5,15] = 1
$data[5,16] = 2
$data[5,17] = "GROUP/DEFAULT"
$data[5,18] = "sword"
$data[5,19] = "null"
$data[5,20] = "null"
$data[5,21] = "END"
$data[6,0] = "Tanuki"
$data[6,1] = "spr_bt_tanuki"
$data[6,2] = "tanooki_1"
$data[6,3] = 30
$data[6,4] = 0
$data[6,5] = 1
$data[6,6] = 1
$data[6,7] = 0
$data[6,8] = 0.75
$data[6,9] = 1
$data[6,10] = 0
$data[6,11] = "Attacker1"
$data[6,12] = "none"
$data[6,13] = 1
$data[6,14] = 0.5
$data[6,15] = 2
$data[6,16] = 1
$data[6,17] = "GROUP/DEFAULT"
$data[6,18] = "lance"
$data[6,19] = "null"
$data[6,20] = "null"
$data[6,21] = "END"
$data[7,0] = "Tanuki2"
$data[7,1] = "spr_bt_tanuki"
$data[7,2] = "tanooki_1"
$data[7,3] = 30
$data[7,4] = 0
$data[7,5] = 2
$data[7,6] = 1
$data[7,7] = 0
$data[7,8] = 0.85
$data[7,9] = 1
$data[7,10] = 25
$data[7,11] = "Attacker1"
$data[7,12] = "none"
$data[7,13] = 1
$data[7,14] = 0.5
$data[7,15] = 2
$data[7,16] = 1
$data[7,17] = "GROUP/DEFAULT"
$data[7,18] = "sword"
$data[7,19] = "null"
$data[7,20] = "null"
$data[7,21] = "lance"
$data[7,22] = "null"
$data[7,23] = "null"
$data[7,24] = "END"
$data[8,0] = "Ijiraq"
$data[8,1] = "spr_bt_ijiraq_placeholder"
$data[8,2] = "ijiraq_1"
$data[8,3] = 20
$data[8,4] = 0
$data[8,5] = 2
$data[8,6] = 1.2
$data[8,7] = 0
$data[8,8] = 0.5
$data[8,9] = 1
$data[8,10] = 0
$data[8,11] = "Attacker1"
$data[8,12] = "none"
$data[8,13] = 0
$data[8,14] = 2
$data[8,15] = -1
$data[8,16] = 1
$data[8,17] = "GROUP/DEFAULT"
$data[8,18] = "sword"
$data[8,19] = "null"
$data[8,20] = "null"
$data[8,21] = "sword"
$data[8,22] = "cryo"
$data[8,23] = "null"
$data[8,24] = "END"
$data[9,0] = "Ijiraq2"
$data[9,1] = "spr_bt_ijiraq_placeholder"
$data[9,2] = "ijiraq_1"
$data[9,3] = 65
$data[9,4] = 0
$data[9,5] = 2
$data[9,6] = 1.2
$data[9,7] = 0
$data[9,8] = 0.6
$data[9,9] = 1
$data[9,10] = 0
$data[9,11] = "Attacker1"
$data[9,12] = "none"
$data[9,13] = 0
$data[9,14] = 2
$data[9,15] = -1
$data[9,16] = 1
$data[9,17] = "GROUP/DEFAULT"
$data[9,18] = "lance"
$data[9,19] = "cryo"
$data[9,20] = "null"
$data[9,21] = "sword"
$data[9,22] = "cryo"
$data[9,23] = "null"
$data[9,24] = "END"
$data[10,0] = "Ijiraq3"
$data[10,1] = "spr_bt_ijiraq_placeholder"
$data[10,2] = "ijiraq_1"
$data[10,3] = 45
$data[10,4] = 0
$data[10,5] = 2
$data[10,6] = 0.75
$data[10,7] = 0
$data[10,8] = 0.5
$data[10,9] = 1
$data[10,10] = 0
$data[10,11] = "Attacker1"
$data[10,12] = "none"
$data[10,13] = 0
$data[10,14] = 1.1000000000000001
$data[10,15] = -1
$data[10,16] = 1
$data[10,17] = "GROUP/DEFAULT"
$data[10,18] = "lance "
$data[10,19] = "cryo"
$data[10,20] = "null"
$data[10,21] = "sword"
$data[10,22] = "cryo"
$data[10,23] = "null"
$data[10,24] = "END"
$data[11,0] = "Doppelganger (BLUE)"
$data[11,1] = "spr_bt_doppelganger_b_placeholder"
$data[11,2] = "doppelganger_1"
$data[11,3] = 100
$data[11,4] = 0
$data[11,5] = 1
$data[11,6] = 1.5
$data[11,7] = 0
$data[11,8] = 0.8
$data[11,9] = 1
$data[11,10] = 0
$data[11,11] = "Doppleganger1"
$data[11,12] = "none"
$data[11,13] = 0.5
$data[11,14] = 2
$data[11,15] = -1
$data[11,16] = 1
$data[11,17] = "GROUP/DEFAULT"
$data[11,18] = "sword"
$data[11,19] = "cryo"
$data[11,20] = "null"
$data[11,21] = "lance"
$data[11,22] = "cryo"
$data[11,23] = "null"
$data[11,24] = "GROUP/TOO_LONG"
$data[11,25] = "hammer"
$data[11,26] = "cryo"
$data[11,27] = "null"
$data[11,28] = "GROUP/SPECIAL"
$data[11,29] = "magic_circle"
$data[11,30] = "null"
$data[11,31] = "null"
$data[11,32] = "END"
$data[12,0] = "Doppelganger (YELLOW)"
$data[12,1] = "spr_bt_doppelganger_y_placeholder"
$data[12,2] = "doppelganger_1"
$data[12,3] = 100
$data[12,4] = 0
$data[12,5] = 2
$data[12,6] = 1.5
$data[12,7] = 0
$data[12,8] = 0.8
$data[12,9] = 1
$data[12,10] = 0
$data[12,11] = "Doppleganger1"
$data[12,12] = "none"
$data[12,13] = 0.5
$data[12,14] = 2
$data[12,15] = -1
$data[12,16] = 1
$data[12,17] = "GROUP/DEFAULT"
$data[12,18] = "sword"
$data[12,19] = "veld"
$data[12,20] = "null"
$data[12,21] = "lance"
$data[12,22] = "veld"
$data[12,23] = "null"
$data[12,24] = "GROUP/TOO_LONG"
$data[12,25] = "hammer"
$data[12,26] = "veld"
$data[12,27] = "null"
$data[12,28] = "GROUP/SPECIAL"
$data[12,29] = "magic_circle"
$data[12,30] = "null"
$data[12,31] = "null"
$data[12,32] = "END"
$data[13,0] = "Doppelganger (RED)"
$data[13,1] = "spr_bt_doppelganger_r_placeholder"
$data[13,2] = "doppelganger_1"
$data[13,3] = 100
$data[13,4] = 0
$data[13,5] = 2
$data[13,6] = 1.5
$data[13,7] = 0
$data[13,8] = 0.8
$data[13,9] = 1
$data[13,10] = 0
$data[13,11] = "Doppleganger1"
$data[13,12] = "none"
$data[13,13] = 0.5
$data[13,14] = 2
$data[13,15] = -1
$data[13,16] = 1
$data[13,17] = "GROUP/DEFAULT"
$data[13,18] = "sword"
$data[13,19] = "agni"
$data[13,20] = "null"
$data[13,21] = "lance"
$data[13,22] = "agni"
$data[13,23] = "null"
$data[13,24] = "GROUP/TOO_LONG"
$data[13,25] = "hammer"
$data[13,26] = "agni"
$data[13,27] = "null"
$data[13,28] = "GROUP/SPECIAL"
$data[13,29] = "magic_circle"
$data[13,30] = "null"
$data[13,31] = "null"
$data[13,32] = "END"
$data[14,0] = "Doppelganger (GRAY)"
$data[14,1] = "spr_bt_doppelganger_g_placeholder"
$data[14,2] = "doppelganger_1"
$data[14,3] = 100
$data[14,4] = 0
$data[14,5] = 2
$data[14,6] = 1.5
$data[14,7] = 0
$data[14,8] = 1.5
$data[14,9] = 1
$data[14,10] = 0
$data[14,11] = "Doppleganger1"
$data[14,12] = "none"
$data[14,13] = 0.5
$data[14,14] = 2
$data[14,15] = -1
$data[14,16] = 1
$data[14,17] = "GROUP/DEFAULT"
$data[14,18] = "hammer"
$data[14,19] = "null"
$data[14,20] = "null"
$data[14,21] = "END"
$data[15,0] = "END"
$ws.Range("A1:AG16").Value = $data
$ws.Range("C15").Select()
